$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames ---
$ws.Range("C1").Value = "startTime"
$ws.Range("E1").Value = "answerKeyword"

# --- Existing-row fixes ---
$ws.Range("D3").Value = "사이코패스"
$ws.Range("C32").Value = 0
$ws.Range("E32").Value = "강철"
$ws.Range("E33").Value = "문호"

# --- Highlight D10 (원피스) in yellow ---
$ws.Range("D10").Interior.Color = 65535

# --- New question rows (39-52) ---
$ws.Range("D39").Value = "철혈의오펜스"
$ws.Range("E39").Value = "철혈 오펜스"

$ws.Range("D40").Value = "알드노아제로"
$ws.Range("E40").Value = "알드노아"

$ws.Range("B41").Value = "VFadUtWFsQk"
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = "가정교사히트맨리본"
$ws.Range("E41").Value = "가히리"

$ws.Range("B42").Value = "0GF5vHEq9LY"
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = "이누야샤"
$ws.Range("E42").Value = "이누야사완결편"

$ws.Range("B43").Value = "IHqlSoqw6mU"
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = "디지몬어드벤쳐"
$ws.Range("E43").Value = "디지몬"

$ws.Range("B44").Value = "nFG3l5zxLdM"
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = "시간을달리는소녀"

$ws.Range("D45").Value = "그날본꽃의이름은우리는아직모른다"
$ws.Range("E45").Value = "아노하나"

$ws.Range("B46").Value = "RXhw8QatPxw"
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = "제로의사역마"

$ws.Range("D47").Value = "신만이아는세계"

$ws.Range("B48").Value = "wLtVF8VUlw8"
$ws.Range("C48").Value = 6
$ws.Range("D48").Value = "러키스타"

$ws.Range("B49").Value = "Ipyi1H1-idg"
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = "인피니트스트라토스"
$ws.Range("E49").Value = "IS"

$ws.Range("B50").Value = "2eJF4o06z58"
$ws.Range("C50").Value = 66
$ws.Range("D50").Value = "엔젤비트"

$ws.Range("B51").Value = "hU8u0BvMzVs"
$ws.Range("C51").Value = 60
$ws.Range("D51").Value = "기교소녀는상처받지않아"

$ws.Range("B52").Value = "3S5Yu9HpE3I"
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = "신세기에반게리온"
$ws.Range("E52").Value = "에반게리온"

# --- View: selection + scroll position ---
$ws.Range("C26").Select()
